$wb = $excel.ActiveWorkbook

$wsStudents = $wb.Worksheets.Item("Студенты")
$wsUniversities = $wb.Worksheets.Item("Университеты")

# --- "Университеты" sheet: fill in previously-blank cells ---
# Row for id=7 (МГМУ): full name was missing, abbreviation already "МГМУ"
$wsUniversities.Range("B4").Value = "Московский Государственный Медицинский Университет"
# Row for id=8 (Тамбовский Университет Медицины): abbreviation was a blank placeholder
$wsUniversities.Range("C5").Value = "ТУМ"

# --- "Студенты" sheet: fill in previously-blank student names ---
$wsStudents.Range("B8").Value = "Витальев В. А."
$wsStudents.Range("B9").Value = "Петров П. А."

# --- Selection / active sheet state ---
# Final state: "Студенты" tab is selected with D1 active; "Университеты" has B13 active.
$null = $wsUniversities.Range("B13").Select()
$null = $wsStudents.Activate()
$null = $wsStudents.Range("D1").Select()
